# "Minor fixes to some of the presentations."
#
# The upstream change swaps the payloads that live behind the two
# SharePoint/"document management" custom XML parts (the part that used to
# hold the <p:properties>/documentManagement stub now holds the
# FormTemplates/contentType stub, and vice-versa) while leaving everything
# else untouched.
#
# Do the swap through the CustomXMLParts COM surface: find each part by its
# namespace, remember its XML, delete both, then re-add them with the
# payloads exchanged so the namespaces (and therefore what each part "is")
# end up swapped.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$propertiesNs = "http://schemas.microsoft.com/office/2006/metadata/properties"
$formsNs      = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"

$propsPart = $null
$formsPart = $null

for ($i = 1; $i -le $parts.Count; $i++) {
    $part = $parts.Item($i)
    if ($part.NamespaceURI -eq $propertiesNs) {
        $propsPart = $part
    } elseif ($part.NamespaceURI -eq $formsNs) {
        $formsPart = $part
    }
}

if ($propsPart -ne $null -and $formsPart -ne $null) {
    $propsXml = $propsPart.XML
    $formsXml = $formsPart.XML

    $propsPart.Delete()
    $formsPart.Delete()

    # Re-create them with the content swapped between the two slots.
    [void]$parts.Add($formsXml)
    [void]$parts.Add($propsXml)
}
